$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from the "At last I added..." paragraph.
#    It will be re-created at the end of the newly appended content below.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Append the new paragraphs at the very end of the document body (after the
#    "Tomorrow I want to have look..." paragraph), with no inherited run/paragraph
#    formatting (plain w:p / w:r, matching the target OOXML) except for the blank
#    separator paragraph, which keeps the en-US language mark used elsewhere in
#    this document.
$endRange = $d.Content
$endRange.Collapse(0)

$newParasXml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' `
  + '<w:p><w:r><w:t>Vandaag layout nog verder mooi gemaakt en custom tab bar icons gemaakt in sketch en deze geupload. Ik kwam</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> er achter dat ze donker van kleur moeten zijn of iig niet blauw anders pakt swift ze niet.</w:t></w:r></w:p>' `
  + '<w:p><w:r><w:t xml:space="preserve">Verder nog meer bugs deruit gehaald zoals het niet weergeven van error messages. Ik heb een alertview toegevoegd aan verschillende viewcontrollers wat ik persoonlijk erg duidelijk vind. </w:t></w:r></w:p>' `
  + '<w:p><w:r><w:t>Verder is de app zo goed als af, heel misschien dat ik morgn nog wat puntjes op de i zet, maar ik ben tevreden eigenlijk.</w:t></w:r>' `
  + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$packageXml = '<?xml version="1.0" standalone="yes"?>' `
  + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
  + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
  + '<pkg:xmlData>' `
  + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
  + '<w:body>' + $newParasXml + '</w:body>' `
  + '</w:document>' `
  + '</pkg:xmlData></pkg:part></pkg:package>'

[void]$endRange.InsertXML($packageXml)

Write-Output "done"
